$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values in row 4 (Persons Receiving Social Package), columns E:K
$ws.Range("E4").Value = 6222
$ws.Range("F4").Value = 6138
$ws.Range("G4").Value = 5987
$ws.Range("H4").Value = 5864
$ws.Range("I4").Value = 6831
$ws.Range("J4").Value = 6994
$ws.Range("K4").Value = 7212

# Update the selected/active cell in the sheet view to A3
$ws.Range("A3").Select()
